# repull data, push all data, mean calculation
# Update the dSF (column F) values for each row per the refreshed pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -2
    4  = 3
    6  = -1
    7  = -2
    8  = 2
    9  = 5
    10 = -7
    11 = 3
    12 = 1
    13 = 5
    14 = 2
    15 = -1
    16 = -1
    17 = 1
    18 = -7
    19 = -3
    21 = -1
    22 = -2
    23 = 1
    25 = -1
    26 = -1
    27 = -2
    28 = 3
    29 = 2
    30 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
